# Swap the order of the two adjacent slides:
#   "MTA traffic every day" (currently slide 4)
#   "Top 10 busiest stations" (currently slide 5)
# After the edit, "Top 10 busiest stations" should come first (position 4)
# and "MTA traffic every day" should follow (position 5).

$p = $ppt.ActivePresentation

# Move the slide currently at position 5 ("Top 10 busiest stations")
# to position 4; PowerPoint shifts "MTA traffic every day" down to
# position 5 automatically.
$p.Slides.Item(5).MoveTo(4)
